$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (AD1:AF1) for the season record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the rest of the header row (bold, bordered,
# centered) by copying the format from an existing header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the season record (Wins/Losses/Ties) for every player row.
$lastRow = 37
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 95   # AD: Wins
    $ws.Cells.Item($row, 31).Value = 67   # AE: Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF: Ties
}

Write-Output "done"
